$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update BLEU score (B11)
$ws.Range("B11").Value = 0.09982210385769184

# Update Code BLEU score (B12) and its note (C12)
$ws.Range("B12").Value = 0.2588900929565529
$ws.Range("C12").Value = "{'codebleu': 0.25889009295655285, 'ngram_match_score': 0.09982210385769184, 'weighted_ngram_match_score': 0.10861581100451317, 'syntax_match_score': 0.5243055555555556, 'dataflow_match_score': 0.3028169014084507}"

# Update Embeddings and Cosine similarity (B13)
$ws.Range("B13").Value = 0.8924435273001929
